# "Generate Report for Handback"
#
# The handback report is regenerated: the file that was just handed back
# (2f846b80-c5f5-4ebe-b969-0e7460431ef6.md) moves to the top data row on
# every sheet, its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and its handback timestamps advance.
# The other two rows (ffff7e577d15... and fffffff47ec0d9...) shift down by
# one row but otherwise keep their values. Hyperlink display text is kept
# in sync with each cell's new text (the underlying r:id/target for a given
# cell position does not change).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @("2f846b80-c5f5-4ebe-b969-0e7460431ef6.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-05-13 11:05:31"),
    @("ffff7e577d15-3540-4630-af39-f5b803a4b64e.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-13 11:03:45"),
    @("fffffff47ec0d9-cbe8-4b94-9dd9-5222ac94e3fe.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-13 11:03:45")
)

for ($i = 0; $i -lt $overviewRows.Count; $i++) {
    $r = $i + 2
    $rowVals = $overviewRows[$i]
    $wsOverview.Cells.Item($r, 1).Value = $rowVals[0]
    $wsOverview.Cells.Item($r, 2).Value = $rowVals[1]
    $wsOverview.Cells.Item($r, 3).Value = $rowVals[2]
    $wsOverview.Cells.Item($r, 4).Value = $rowVals[3]
}

$overviewLinks = @()
foreach ($hl in $wsOverview.Hyperlinks) { $overviewLinks += $hl }
# Order: A2, A3, A4
$overviewLinks[0].TextToDisplay = $overviewRows[0][0]
$overviewLinks[1].TextToDisplay = $overviewRows[1][0]
$overviewLinks[2].TextToDisplay = $overviewRows[2][0]

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhRows = @(
    @("2f846b80-c5f5-4ebe-b969-0e7460431ef6.md", ".md", "Handed back: in sync with en-US", `
      "2f846b80-c5f5-4ebe-b969-0e7460431ef6.87eea981de8f071153fd104ea2bec243947365cf.zh-cn.xlf", `
      "2016-03-13 11:05:27", "2f846b80-c5f5-4ebe-b969-0e7460431ef6.md", `
      "2f846b80-c5f5-4ebe-b969-0e7460431ef6.87eea981de8f071153fd104ea2bec243947365cf.zh-cn.xlf", `
      "2016-03-13 11:05:44", "Include"),
    @("ffff7e577d15-3540-4630-af39-f5b803a4b64e.md", ".md", "Handed back: in sync with en-US", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf", `
      "2016-03-13 11:03:41", "2246d8cb-028b-463f-8a0a-0d8d45762021.md", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf", `
      "2016-03-13 11:03:58", "Include"),
    @("fffffff47ec0d9-cbe8-4b94-9dd9-5222ac94e3fe.md", ".md", "Handed back: in sync with en-US", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf", `
      "2016-03-13 11:03:41", "2246d8cb-028b-463f-8a0a-0d8d45762021.md", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf", `
      "2016-03-13 11:03:58", "Include")
)

for ($i = 0; $i -lt $zhRows.Count; $i++) {
    $r = $i + 2
    $rowVals = $zhRows[$i]
    $wsZh.Cells.Item($r, 1).Value = $rowVals[0]
    $wsZh.Cells.Item($r, 2).Value = $rowVals[1]
    $wsZh.Cells.Item($r, 3).Value = $rowVals[2]
    $wsZh.Cells.Item($r, 4).Value = $rowVals[3]
    $wsZh.Cells.Item($r, 5).Value = $rowVals[4]
    $wsZh.Cells.Item($r, 6).Value = $rowVals[5]
    $wsZh.Cells.Item($r, 7).Value = $rowVals[6]
    $wsZh.Cells.Item($r, 8).Value = $rowVals[7]
    $wsZh.Cells.Item($r, 9).Value = $rowVals[8]
}

$zhLinks = @()
foreach ($hl in $wsZh.Hyperlinks) { $zhLinks += $hl }
# Order per row: A, B, D, F, G  (5 links per row x 3 rows = 15)
for ($i = 0; $i -lt $zhRows.Count; $i++) {
    $rowVals = $zhRows[$i]
    $base = $i * 5
    $zhLinks[$base + 0].TextToDisplay = $rowVals[0]
    $zhLinks[$base + 1].TextToDisplay = $rowVals[1]
    $zhLinks[$base + 2].TextToDisplay = $rowVals[3]
    $zhLinks[$base + 3].TextToDisplay = $rowVals[5]
    $zhLinks[$base + 4].TextToDisplay = $rowVals[6]
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deRows = @(
    @("2f846b80-c5f5-4ebe-b969-0e7460431ef6.md", ".md", "Handed back: in sync with en-US", `
      "2f846b80-c5f5-4ebe-b969-0e7460431ef6.87eea981de8f071153fd104ea2bec243947365cf.de-de.xlf", `
      "2016-03-13 11:05:31", "2f846b80-c5f5-4ebe-b969-0e7460431ef6.md", `
      "2f846b80-c5f5-4ebe-b969-0e7460431ef6.87eea981de8f071153fd104ea2bec243947365cf.de-de.xlf", `
      "2016-03-13 11:05:50", "Include"),
    @("ffff7e577d15-3540-4630-af39-f5b803a4b64e.md", ".md", "Handed back: in sync with en-US", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf", `
      "2016-03-13 11:03:45", "2246d8cb-028b-463f-8a0a-0d8d45762021.md", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf", `
      "2016-03-13 11:04:05", "Include"),
    @("fffffff47ec0d9-cbe8-4b94-9dd9-5222ac94e3fe.md", ".md", "Handed back: in sync with en-US", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf", `
      "2016-03-13 11:03:45", "2246d8cb-028b-463f-8a0a-0d8d45762021.md", `
      "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf", `
      "2016-03-13 11:04:05", "Include")
)

for ($i = 0; $i -lt $deRows.Count; $i++) {
    $r = $i + 2
    $rowVals = $deRows[$i]
    $wsDe.Cells.Item($r, 1).Value = $rowVals[0]
    $wsDe.Cells.Item($r, 2).Value = $rowVals[1]
    $wsDe.Cells.Item($r, 3).Value = $rowVals[2]
    $wsDe.Cells.Item($r, 4).Value = $rowVals[3]
    $wsDe.Cells.Item($r, 5).Value = $rowVals[4]
    $wsDe.Cells.Item($r, 6).Value = $rowVals[5]
    $wsDe.Cells.Item($r, 7).Value = $rowVals[6]
    $wsDe.Cells.Item($r, 8).Value = $rowVals[7]
    $wsDe.Cells.Item($r, 9).Value = $rowVals[8]
}

$deLinks = @()
foreach ($hl in $wsDe.Hyperlinks) { $deLinks += $hl }
for ($i = 0; $i -lt $deRows.Count; $i++) {
    $rowVals = $deRows[$i]
    $base = $i * 5
    $deLinks[$base + 0].TextToDisplay = $rowVals[0]
    $deLinks[$base + 1].TextToDisplay = $rowVals[1]
    $deLinks[$base + 2].TextToDisplay = $rowVals[3]
    $deLinks[$base + 3].TextToDisplay = $rowVals[5]
    $deLinks[$base + 4].TextToDisplay = $rowVals[6]
}

Write-Output "Report regenerated for handback."
